# Apply scheduled-runner price/profit updates to Famfrit_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3308.8572
$ws.Range("I38").Value = 2618.9333
$ws.Range("J38").Value = 5033.6665
$ws.Range("K38").Value = 7856.7999
$ws.Range("L38").Value = 15100.9995
$ws.Range("M38").Value = -7484.7999
$ws.Range("N38").Value = -15844.9995
$ws.Range("H70").Value = 3450
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3450
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10350
$ws.Range("N70").Value = -10890
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 3450
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3450
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10350
$ws.Range("N73").Value = -12222
$ws.Range("M73").ClearContents()
$ws.Range("H74").Value = 4642
$ws.Range("I74").Value = 2494.6
$ws.Range("J74").Value = 5715.7
$ws.Range("K74").Value = 2494.6
$ws.Range("L74").Value = 5715.7
$ws.Range("M74").Value = -1558.6
$ws.Range("N74").Value = -7587.7
$ws.Range("H77").Value = 4642
$ws.Range("I77").Value = 2494.6
$ws.Range("J77").Value = 5715.7
$ws.Range("K77").Value = 12473
$ws.Range("L77").Value = 28578.5
$ws.Range("M77").Value = -7793
$ws.Range("N77").Value = -37938.5
$ws.Range("H132").Value = 3137.4285
$ws.Range("I132").Value = 2731.4211
$ws.Range("J132").Value = 6994.5
$ws.Range("K132").Value = 8194.263300000001
$ws.Range("L132").Value = 20983.5
$ws.Range("M132").Value = -5664.263300000001
$ws.Range("N132").Value = -26043.5
$ws.Range("H136").Value = 76365.57000000001
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 76365.57000000001
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 76365.57000000001
$ws.Range("N136").Value = -86565.57000000001
$ws.Range("H137").Value = 2125.1936
$ws.Range("I137").Value = 2138.8096
$ws.Range("J137").Value = 2096.6
$ws.Range("K137").Value = 6416.4288
$ws.Range("L137").Value = 6289.799999999999
$ws.Range("M137").Value = -3866.4288
$ws.Range("N137").Value = -11389.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9430.429
$ws.Range("I32").Value = 6705.7173
$ws.Range("J32").Value = 21964.1
$ws.Range("K32").Value = 6705.7173
$ws.Range("L32").Value = 21964.1
$ws.Range("M32").Value = -6418.7173
$ws.Range("N32").Value = -22538.1
$ws.Range("H110").Value = 32045.75
$ws.Range("I110").Value = 40951
$ws.Range("J110").Value = 5330
$ws.Range("K110").Value = 40951
$ws.Range("L110").Value = 5330
$ws.Range("M110").Value = -38906
$ws.Range("N110").Value = -9420
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 108
$ws.Range("I5").Value = 108
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 108
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 5
$ws.Range("N5").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H75").Value = 25333.334
$ws.Range("I75").Value = 11000
$ws.Range("J75").Value = 54000
$ws.Range("K75").Value = 11000
$ws.Range("L75").Value = 54000
$ws.Range("M75").Value = -10064
$ws.Range("N75").Value = -55872
$ws.Range("H78").Value = 25333.334
$ws.Range("I78").Value = 11000
$ws.Range("J78").Value = 54000
$ws.Range("K78").Value = 33000
$ws.Range("L78").Value = 162000
$ws.Range("M78").Value = -28320
$ws.Range("N78").Value = -171360
$ws.Range("H80").Value = 1118
$ws.Range("I80").Value = 776.875
$ws.Range("J80").Value = 1421.2222
$ws.Range("K80").Value = 776.875
$ws.Range("L80").Value = 1421.2222
$ws.Range("M80").Value = 221.125
$ws.Range("N80").Value = -3417.2222
$ws.Range("H83").Value = 1118
$ws.Range("I83").Value = 776.875
$ws.Range("J83").Value = 1421.2222
$ws.Range("K83").Value = 3884.375
$ws.Range("L83").Value = 7106.111
$ws.Range("M83").Value = 1107.625
$ws.Range("N83").Value = -17090.111
$ws.Range("H87").Value = 92750
$ws.Range("I87").Value = 95500
$ws.Range("J87").Value = 90000
$ws.Range("K87").Value = 95500
$ws.Range("L87").Value = 90000
$ws.Range("M87").Value = -94252
$ws.Range("N87").Value = -92496
$ws.Range("H90").Value = 92750
$ws.Range("I90").Value = 95500
$ws.Range("J90").Value = 90000
$ws.Range("K90").Value = 286500
$ws.Range("L90").Value = 270000
$ws.Range("M90").Value = -280260
$ws.Range("N90").Value = -282480
$ws.Range("H94").Value = 1066.5
$ws.Range("I94").Value = 264.77777
$ws.Range("J94").Value = 1722.4546
$ws.Range("K94").Value = 264.77777
$ws.Range("L94").Value = 1722.4546
$ws.Range("M94").Value = 186.22223
$ws.Range("N94").Value = -2624.4546
$ws.Range("H132").Value = 290000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 290000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 290000
$ws.Range("N132").Value = -300120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 52682.85
$ws.Range("I132").Value = 73415.25
$ws.Range("J132").Value = 4307.25
$ws.Range("K132").Value = 220245.75
$ws.Range("L132").Value = 12921.75
$ws.Range("M132").Value = -217715.75
$ws.Range("N132").Value = -17981.75
$ws.Range("H134").Value = 2098.625
$ws.Range("I134").Value = 1921.45
$ws.Range("J134").Value = 2984.5
$ws.Range("K134").Value = 5764.35
$ws.Range("L134").Value = 8953.5
$ws.Range("M134").Value = -3229.35
$ws.Range("N134").Value = -14023.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 6064.9
$ws.Range("I50").Value = 300
$ws.Range("J50").Value = 6705.4443
$ws.Range("K50").Value = 900
$ws.Range("L50").Value = 20116.3329
$ws.Range("M50").Value = -419
$ws.Range("N50").Value = -21078.3329
$ws.Range("H53").Value = 6064.9
$ws.Range("I53").Value = 300
$ws.Range("J53").Value = 6705.4443
$ws.Range("K53").Value = 900
$ws.Range("L53").Value = 20116.3329
$ws.Range("M53").Value = -419
$ws.Range("N53").Value = -21078.3329
$ws.Range("H87").Value = 1250.5
$ws.Range("I87").Value = 1250.5
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 3751.5
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -2503.5
$ws.Range("H90").Value = 1250.5
$ws.Range("I90").Value = 1250.5
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 11254.5
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -5014.5
$ws.Range("H98").Value = 2017.8572
$ws.Range("I98").Value = 225
$ws.Range("J98").Value = 2316.6667
$ws.Range("K98").Value = 675
$ws.Range("L98").Value = 6950.000100000001
$ws.Range("M98").Value = 823
$ws.Range("N98").Value = -9946.000100000001
$ws.Range("H121").Value = 6482012
$ws.Range("I121").Value = 324.7143
$ws.Range("J121").Value = 15556375
$ws.Range("K121").Value = 974.1428999999999
$ws.Range("L121").Value = 46669125
$ws.Range("M121").Value = 335.8571000000001
$ws.Range("N121").Value = -46671745
$ws.Range("H127").Value = 3739.6
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 3739.6
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 11218.8
$ws.Range("N127").Value = -21138.8
$ws.Range("H131").Value = 1715.129
$ws.Range("I131").Value = 1650
$ws.Range("J131").Value = 1722.1072
$ws.Range("K131").Value = 4950
$ws.Range("L131").Value = 5166.321599999999
$ws.Range("M131").Value = 90
$ws.Range("N131").Value = -15246.3216
$ws.Range("H134").Value = 1275.7894
$ws.Range("I134").Value = 1275.7894
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3827.3682
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1242.6318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 164470.67
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 164470.67
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 164470.67
$ws.Range("N95").Value = -169962.67
$ws.Range("H131").Value = 89500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 89500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 89500
$ws.Range("N131").Value = -99580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 651.3913
$ws.Range("I55").Value = 408.8
$ws.Range("J55").Value = 1106.25
$ws.Range("K55").Value = 408.8
$ws.Range("L55").Value = 1106.25
$ws.Range("M55").Value = -235.8
$ws.Range("N55").Value = -1452.25
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 3116.7727
$ws.Range("I82").Value = 2269.7273
$ws.Range("J82").Value = 3963.818
$ws.Range("K82").Value = 2269.7273
$ws.Range("L82").Value = 3963.818
$ws.Range("M82").Value = -1908.7273
$ws.Range("N82").Value = -4685.818
$ws.Range("H85").Value = 3116.7727
$ws.Range("I85").Value = 2269.7273
$ws.Range("J85").Value = 3963.818
$ws.Range("K85").Value = 2269.7273
$ws.Range("L85").Value = 3963.818
$ws.Range("M85").Value = -1021.7273
$ws.Range("N85").Value = -6459.818
$ws.Range("H93").Value = 2202.0588
$ws.Range("I93").Value = 1137.6
$ws.Range("J93").Value = 2645.5833
$ws.Range("K93").Value = 1137.6
$ws.Range("L93").Value = 2645.5833
$ws.Range("M93").Value = 110.4000000000001
$ws.Range("N93").Value = -5141.5833
$ws.Range("H96").Value = 42598.5
$ws.Range("I96").Value = 30000
$ws.Range("J96").Value = 55197
$ws.Range("K96").Value = 30000
$ws.Range("L96").Value = 55197
$ws.Range("M96").Value = -27254
$ws.Range("N96").Value = -60689
$ws.Range("H132").Value = 19109.242
$ws.Range("I132").Value = 20210.715
$ws.Range("J132").Value = 18297.63
$ws.Range("K132").Value = 60632.145
$ws.Range("L132").Value = 54892.89
$ws.Range("M132").Value = -58102.145
$ws.Range("N132").Value = -59952.89

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 336.875
$ws.Range("I107").Value = 336.875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1010.625
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 909.375
$ws.Range("H132").Value = 3672.111
$ws.Range("I132").Value = 4435.5713
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 13306.7139
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -10776.7139
$ws.Range("N132").Value = -8060

